# Assignment 2 finalization: fill in the Developer name and the
# Method Inputs / Method Call / Expected Result columns of the unit
# test plan for the InvestmentAccount class, then leave the selection
# on G9 (matching the author's last saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name (row 4's label is in B4/C4 "Class Name:"; row 3 is "Developer:")
$ws.Range("C3").Value = "Md Apurba Khan"

# Test case 1 (row 7) - __init__ / valid inputs
$ws.Range("E7").Value = "Valid account_number, client_number, balance, date_created, management_fee"
$ws.Range("F7").Value = "InvestmentAccount(4001, 1002, 1000, date(2010, 1, 1), 2.5)"
$ws.Range("G7").Value = "Instance is created successfully with correct attributes."

# Test case 2 (row 8) - __init__ / invalid management fee type
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'InvestmentAccount(4001, 1002, 1000, date(2015, 6, 1), "invalid")'
$ws.Range("G8").Value = "Management fee defaults to 2.55."

# Test case 3 (row 9) - get_service_charges / date created more than 10 years ago
$ws.Range("E9").Value = "self._date_created = date(2000, 1, 1)"
$ws.Range("F9").Value = "get_service_charges()"
$ws.Range("G9").Value = "Returns base service charge (e.g., `$0.50), management fee is waived."

# Test case 4 (row 10) - get_service_charges / date created within last 10 years
$ws.Range("E10").Value = "self._date_created = date.today() - timedelta(days=10*365.25)"
$ws.Range("F10").Value = "get_service_charges()"
$ws.Range("G10").Value = "Returns base service charge + management fee."

# Test case 5 (row 11) - get_service_charges / date created exactly 10 years ago
$ws.Range("E11").Value = "self._date_created = date(2020, 1, 1)"
$ws.Range("F11").Value = "get_service_charges()"
$ws.Range("G11").Value = "Returns base service charge + management fee."

# Test case 6 (row 12) - __str__ / waived management fee display
$ws.Range("E12").Value = "self._date_created = date(2000, 1, 1)"
$ws.Range("F12").Value = "str(investment_account)"
$ws.Range("G12").Value = "Account details show service charge as `$0.50, no management fee applied."

# Test case 7 (row 13) - __str__ / applied management fee display
$ws.Range("E13").Value = "self._date_created = date(2020, 1, 1)"
$ws.Range("F13").Value = "str(investment_account)"
$ws.Range("G13").Value = "Account details show applied management fee in addition to base charge."

# Restore cursor/selection to where the author left it when they saved.
$null = $ws.Range("G9").Select()
